# "improve second example analysis"
# Insert a "source" column (E) into the gresham mass profiles sheet, recompute
# the baseline rows against a 2.04722-based conversion factor, and add four
# new scenario rows (compost_1000t, compost_6pt5pct, reduce_fw_05pct,
# reduce_fw_10pct).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gresham mass profiles")

# --- shift tons/miles/notes one column right, insert new "source" column E ---
$ws.Columns("E").Insert()

# --- column widths ---
$ws.Columns("A").ColumnWidth = 16.6
$ws.Columns("E").ColumnWidth = 13.33

# --- header row ---
$ws.Range("E1").Value = "umbDisp"

# --- row 2: baseline / FoodWaste / landfilling ---
$ws.Range("A2").Value = "baseline"
$ws.Range("E2").Value = "disposal"
$ws.Range("F2").Formula = "=ROUND(SUM('metro mass profile 2018'!E2:E6)*greshamMetroRatio/2.04722,0)"
$ws.Range("G2").Value = 178
$ws.Range("H2").Value = "based on Metro's total and Gresham's population as a portion of Metro's, adjusted to 9000 tons yard debris standard"

# --- row 3: baseline / YardDebris / composting ---
$ws.Range("A3").Value = "baseline"
$ws.Range("E3").Value = "recovery"
$ws.Range("F3").Formula = "=ROUND(SUM('metro mass profile 2018'!E7:E11)*greshamMetroRatio/2.047222,0)"
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = "based on Metro's total and Gresham's population as a portion of Metro's, adjusted to 9000 tons yard debris standard"

# --- fill in the new scenario names first (column A), top to bottom for the
#     compost rows, then the 10pct-before-5pct ordering the author happened
#     to use, which is what drives shared-string insertion order ---
$ws.Range("A4").Value = "compost_1000t"
$ws.Range("A5").Value = "compost_1000t"
$ws.Range("A6").Value = "compost_1000t"
$ws.Range("A7").Value = "compost_6pt5pct"
$ws.Range("A8").Value = "compost_6pt5pct"
$ws.Range("A9").Value = "compost_6pt5pct"
$ws.Range("A12").Value = "reduce_fw_10pct"
$ws.Range("A13").Value = "reduce_fw_10pct"
$ws.Range("A10").Value = "reduce_fw_05pct"
$ws.Range("A11").Value = "reduce_fw_05pct"

# --- row 4: compost_1000t / FoodWaste / composting ---
$ws.Range("B4").Value = "Gresham"
$ws.Range("C4").Value = "FoodWaste"
$ws.Range("D4").Value = "composting"
$ws.Range("E4").Value = "recovery"
$ws.Range("F4").Value = 1000
$ws.Range("G4").Value = 78

# --- row 5: compost_1000t / FoodWaste / landfilling ---
$ws.Range("B5").Value = "Gresham"
$ws.Range("C5").Value = "FoodWaste"
$ws.Range("D5").Value = "landfilling"
$ws.Range("E5").Value = "disposal"
$ws.Range("F5").Formula = "=F2-F4"
$ws.Range("G5").Value = 178

# --- row 6: compost_1000t / YardDebris / composting ---
$ws.Range("B6").Value = "Gresham"
$ws.Range("C6").Value = "YardDebris"
$ws.Range("D6").Value = "composting"
$ws.Range("E6").Value = "recovery"
$ws.Range("F6").Value = 9000
$ws.Range("G6").Value = 78

# --- row 7: compost_6pt5pct / FoodWaste / composting ---
$ws.Range("B7").Value = "Gresham"
$ws.Range("C7").Value = "FoodWaste"
$ws.Range("D7").Value = "composting"
$ws.Range("E7").Value = "recovery"
$ws.Range("F7").Formula = "=F9*0.065"
$ws.Range("G7").Value = 78

# --- row 8: compost_6pt5pct / FoodWaste / landfilling ---
$ws.Range("B8").Value = "Gresham"
$ws.Range("C8").Value = "FoodWaste"
$ws.Range("D8").Value = "landfilling"
$ws.Range("E8").Value = "disposal"
$ws.Range("F8").Formula = "=F2-F7"
$ws.Range("G8").Value = 178

# --- row 9: compost_6pt5pct / YardDebris / composting ---
$ws.Range("B9").Value = "Gresham"
$ws.Range("C9").Value = "YardDebris"
$ws.Range("D9").Value = "composting"
$ws.Range("E9").Value = "recovery"
$ws.Range("F9").Value = 9000
$ws.Range("G9").Value = 78

# --- row 10: reduce_fw_05pct / FoodWaste / landfilling (new) ---
$ws.Range("B10").Value = "Gresham"
$ws.Range("C10").Value = "FoodWaste"
$ws.Range("D10").Value = "landfilling"
$ws.Range("E10").Value = "disposal"
$ws.Range("F10").Formula = "=F2*0.95"
$ws.Range("G10").Value = 178

# --- row 11: reduce_fw_05pct / YardDebris / composting (new) ---
$ws.Range("B11").Value = "Gresham"
$ws.Range("C11").Value = "YardDebris"
$ws.Range("D11").Value = "composting"
$ws.Range("E11").Value = "recovery"
$ws.Range("F11").Formula = "=F3"
$ws.Range("G11").Value = 5

# --- row 12: reduce_fw_10pct / FoodWaste / landfilling (new) ---
$ws.Range("B12").Value = "Gresham"
$ws.Range("C12").Value = "FoodWaste"
$ws.Range("D12").Value = "landfilling"
$ws.Range("E12").Value = "disposal"
$ws.Range("F12").Formula = "=F2*0.9"
$ws.Range("G12").Value = 178

# --- row 13: reduce_fw_10pct / YardDebris / composting (new) ---
$ws.Range("B13").Value = "Gresham"
$ws.Range("C13").Value = "YardDebris"
$ws.Range("D13").Value = "composting"
$ws.Range("E13").Value = "recovery"
$ws.Range("F13").Formula = "=F3"
$ws.Range("G13").Value = 5

# --- notes column, filled after all the new rows' other cells, matching the
#     order new note strings were introduced ---
$ws.Range("H4").Value = "based on projection that 1000 tons of food waste could be composted"
$ws.Range("H5").Value = "based on projection that 1000 tons of food waste could be composted"
$ws.Range("H6").Value = "based on projection that 1000 tons of food waste could be composted"
$ws.Range("H7").Value = "based on idea that food waste is 6.5% of mixed yard debris/food waste (see metro waste study)"
$ws.Range("H8").Value = "based on idea that food waste is 6.5% of mixed yard debris/food waste (see metro waste study)"
$ws.Range("H9").Value = "based on idea that food waste is 6.5% of mixed yard debris/food waste (see metro waste study)"
$ws.Range("H10").Value = "baseline, but with 5 percent less food waste"
$ws.Range("H11").Value = "baseline, but with 5 percent less food waste"
$ws.Range("H12").Value = "baseline, but with 10 percent less food waste"
$ws.Range("H13").Value = "baseline, but with 10 percent less food waste"

# --- selection to match the saved file ---
$ws.Range("F7").Select()
